$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 30: drop "R58" from the (rich-text) resistor list in C30 ---------
# Before: R10,R13,R21,R23,R39,R50,R51,R58,R59,R62,R64
#   runs:  [default]          [green]R39[/]  [black],R50,R51,R58,[/]  [green]R59[/]  [black],R62,[/]  [red]R64[/]
# After : R10,R13,R21,R23,R39,R50,R51,R59,R62,R64
#   runs:  [default]          [green]R39[/]  [black],R50,R51,[/]      [green]R59[/]  [black],R62,[/]  [red]R64[/]
$newC30 = "R10,R13,R21,R23,R39,R50,R51,R59,R62,R64"
$c30 = $ws.Range("C30")
$c30.Value2 = $newC30

# Re-apply the per-run colouring that existed before the edit (everything
# else keeps the cell's default black font, same as the original "R10,R13,
# R21,R23," lead-in run).
$i39 = $newC30.IndexOf("R39") + 1
$c30.Characters($i39, 3).Font.Color = 0x50B000   # green FF00B050

$iMid = $newC30.IndexOf(",R50,R51,") + 1
$c30.Characters($iMid, 9).Font.Color = 0         # black FF000000

$i59 = $newC30.LastIndexOf("R59") + 1
$c30.Characters($i59, 3).Font.Color = 0x50B000   # green FF00B050

$iTail = $newC30.LastIndexOf(",R62,") + 1
$c30.Characters($iTail, 5).Font.Color = 0        # black FF000000

$i64 = $newC30.LastIndexOf("R64") + 1
$c30.Characters($i64, 3).Font.Color = 0xFF       # red FF0000

# --- Row 35: drop "R56" from the (plain-text) resistor list in C35 --------
# Before: R11,R14,R37,R48,R49,R56
# After : R11,R14,R37,R48,R49
$ws.Range("C35").Value2 = "R11,R14,R37,R48,R49"

# --- Restore the view state (scrolled down, selection moved to C38) -------
[void]$ws.Range("C38").Select()
